# Refresh the cryptocurrency price/volume snapshot (cols D = Price, E = Volume(1h))
# on Sheet1, mirroring the GitHub Actions scheduled data refresh.
#
# Several "Price" values (col D) look like plain numbers to Excel's parser
# (e.g. "403.13", "1.00"), which would otherwise get auto-coerced to a
# number and lose the exact text formatting. To keep them as literal text
# -- matching the workbook's inlineStr cells -- we briefly force the cell
# to Text format before assigning, then restore the default "Normal" style
# so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.603.52"
$ws.Range("E2").Value = "  +5.59%  "

$ws.Range("D3").Value = "3.188.28"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "403.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.28%  "

$ws.Range("E7").Value = "  +1.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.11%  "

$ws.Range("E11").Value = "  +1.71%  "

$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").Value = "3.679.26"
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("E16").Value = "  +8.71%  "

$ws.Range("D17").Value = "3.188.57"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("D19").Value = "54.493.08"
$ws.Range("E19").Value = "  +5.06%  "

$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.94%  "

$ws.Range("D22").Value = "0.0₃0996"
$ws.Range("E22").Value = "  +2.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("E25").Value = "  +4.64%  "

$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.20%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +4.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0500"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.41%  "

$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("E37").Value = "  +8.64%  "

$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("E39").Value = "  +9.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.01%  "

$ws.Range("E42").Value = "  -1.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "

$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.04%  "

$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("D49").Value = "2.089.30"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("E50").Value = "  +7.73%  "

$ws.Range("E51").Value = "  +10.80%  "
